# Update cosinor analysis results (sawtooth_05) after re-running CircadiPy simulations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 22.70000000000011
$ws.Range("G2").Value = [double]"3.648916705456529e-06"
$ws.Range("H2").Value = [double]"2.101619151221488e-05"
$ws.Range("I2").Value = 0.1397219266286147
$ws.Range("K2").Value = 4.846680087310672
$ws.Range("L2").Value = "[2.3087930397586565, 7.384567134862687]"
$ws.Range("M2").Value = 0.0001983204816922157
$ws.Range("N2").Value = 0.0001983204816922157
$ws.Range("O2").Value = -0.9811580659860013
$ws.Range("P2").Value = "[-1.509473947670771, -0.4528421843012316]"
$ws.Range("Q2").Value = 0.0002943179237711036
$ws.Range("R2").Value = 0.0003049156030161537
$ws.Range("S2").Value = 11.54651635097832
$ws.Range("T2").Value = "[10.220423141405096, 12.872609560551542]"
$ws.Range("W2").Value = 3.544744744744762
$ws.Range("X2").Value = 1.636036036036045
$ws.Range("Y2").Value = 5.453453453453479

# ---- Row 3 ----
$ws.Range("E3").Value = 23.23000000000019
$ws.Range("G3").Value = [double]"7.513966182592924e-07"
$ws.Range("H3").Value = [double]"1.409519809018415e-05"
$ws.Range("K3").Value = 5.836211037737688
$ws.Range("L3").Value = "[2.898289876973704, 8.774132198501672]"
$ws.Range("M3").Value = 0.000112885007220509
$ws.Range("N3").Value = 0.0001983204816922157
$ws.Range("O3").Value = 0.8868159442565782
$ws.Range("P3").Value = "[0.40881586082750143, 1.364816027685655]"
$ws.Range("Q3").Value = 0.0003049156030161537
$ws.Range("R3").Value = 0.0003049156030161537
$ws.Range("S3").Value = 11.62868285882105
$ws.Range("T3").Value = "[10.120680820031009, 13.136684897611087]"
$ws.Range("W3").Value = 19.95129129129146
$ws.Range("X3").Value = 18.1840440440442
$ws.Range("Y3").Value = 21.71853853853872
